$wb = $excel.ActiveWorkbook

# Rename the existing sheet from "Sheet1" to "factors"
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "factors"

# Insert a new sheet "R2M" right after "factors" and make it active
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "R2M"

# Fill in the small R2S concordance table (column-major entry order so
# shared strings land in the same slot order as the source edit).
$ws2.Range("A1").Value = "aggregator"
$ws2.Range("A2").Value = "industry"
$ws2.Range("A3").Value = "product"
$ws2.Range("B1").Value = "date"

$ws2.Range("A1:B1").Font.Bold = $true

$ws2.Range("B2").Value = 20200421
$ws2.Range("B3").Value = 20200421

# Match the header column's fitted width
$ws2.Columns.Item(1).ColumnWidth = 8.8

# Page setup to mirror the workbook's other sheet
$ws2.PageSetup.PaperSize = 9
$ws2.PageSetup.Orientation = 1

# Leave the selection where the author left off, on the new sheet
$ws2.Range("B4").Select() | Out-Null
